# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45679 (2025-01-22) to 45680 (2025-01-23).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 36; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45679) {
        $cell.Value2 = 45680
    }
}
